# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks numeric: force Text format first so Excel
# does not silently convert the literal (with trailing zeros / exact
# decimal digits) into a binary double, then strip the format again so
# no stray style index is left behind.
$numericTextCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D14", "D17", "D19", "D22", "D23", "D25", "D27", "D30", "D31", "D32", "D33", "D35", "D36", "D39", "D40", "D42", "D43", "D47", "D48", "D50", "D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "44.711.67"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "2.421.72"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "316.28"
$ws.Range("E5").Value = "  +4.52%  "
$ws.Range("D6").Value = "101.02"
$ws.Range("E6").Value = "  +6.25%  "
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  +11.63%  "
$ws.Range("D10").Value = "35.31"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("D11").Value = "0.0798"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "18.68"
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "6.90"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").Value = "2.798.46"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "2.423.42"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").Value = "0.830"
$ws.Range("E17").Value = "  +4.61%  "
$ws.Range("D18").Value = "44.515.66"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("D19").Value = "12.31"
$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "0.0₃0915"
$ws.Range("E21").Value = "  +3.56%  "
$ws.Range("D22").Value = "68.67"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "242.17"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "25.18"
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").Value = "33.47"
$ws.Range("E30").Value = "  +4.04%  "
$ws.Range("D31").Value = "48.46"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "0.126"
$ws.Range("E32").Value = "  +17.86%  "
$ws.Range("D33").Value = "19.44"
$ws.Range("E33").Value = "  +11.20%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0769"
$ws.Range("E36").Value = "  +6.23%  "
$ws.Range("E37").Value = "  +3.68%  "
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").Value = "2.85"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "124.93"
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("D43").Value = "21.08"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("D45").Value = "1.938.96"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "2.93"
$ws.Range("E47").Value = "  +8.14%  "
$ws.Range("D48").Value = "9.23"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  +15.22%  "
$ws.Range("D50").Value = "75.69"
$ws.Range("E50").Value = "  +6.08%  "
$ws.Range("D51").Value = "53.89"
$ws.Range("E51").Value = "  +5.60%  "

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).ClearFormats()
}

